$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value2 = 1682.8572  # H98: 2825 -> 1682.8572
$ws.Cells.Item(98, 9).Value2 = 1196  # I98: 2000 -> 1196
$ws.Cells.Item(98, 10).Value2 = 2900  # J98: 3100 -> 2900
$ws.Cells.Item(98, 11).Value2 = 1196  # K98: 2000 -> 1196
$ws.Cells.Item(98, 12).Value2 = 2900  # L98: 3100 -> 2900
$ws.Cells.Item(98, 13).Value2 = 302  # M98: -502 -> 302
$ws.Cells.Item(98, 14).Value2 = -5896  # N98: -6096 -> -5896
$ws.Cells.Item(113, 8).Value2 = 40003932  # H113: 37040828 -> 40003932
$ws.Cells.Item(113, 9).Value2 = 76926360  # I113: 71431770 -> 76926360
$ws.Cells.Item(113, 10).Value2 = 4633.1665  # J113: 4427.5386 -> 4633.1665
$ws.Cells.Item(113, 11).Value2 = 76926360  # K113: 71431770 -> 76926360
$ws.Cells.Item(113, 12).Value2 = 4633.1665  # L113: 4427.5386 -> 4633.1665
$ws.Cells.Item(113, 13).Value2 = -76923106  # M113: -71428516 -> -76923106
$ws.Cells.Item(113, 14).Value2 = -11141.1665  # N113: -10935.5386 -> -11141.1665
$ws.Cells.Item(122, 8).Value2 = 1682.8572  # H122: 2825 -> 1682.8572
$ws.Cells.Item(122, 9).Value2 = 1196  # I122: 2000 -> 1196
$ws.Cells.Item(122, 10).Value2 = 2900  # J122: 3100 -> 2900
$ws.Cells.Item(122, 11).Value2 = 3588  # K122: 6000 -> 3588
$ws.Cells.Item(122, 12).Value2 = 8700  # L122: 9300 -> 8700
$ws.Cells.Item(122, 13).Value2 = -1138  # M122: -3550 -> -1138
$ws.Cells.Item(122, 14).Value2 = -13600  # N122: -14200 -> -13600
$ws.Cells.Item(132, 8).Value2 = 1668.0377  # H132: 1906.8914 -> 1668.0377
$ws.Cells.Item(132, 9).Value2 = 1792.3673  # I132: 2028.7675 -> 1792.3673
$ws.Cells.Item(132, 10).Value2 = 145  # J132: 160 -> 145
$ws.Cells.Item(132, 11).Value2 = 5377.1019  # K132: 6086.3025 -> 5377.1019
$ws.Cells.Item(132, 12).Value2 = 435  # L132: 480 -> 435
$ws.Cells.Item(132, 13).Value2 = -2847.1019  # M132: -3556.3025 -> -2847.1019
$ws.Cells.Item(132, 14).Value2 = -5495  # N132: -5540 -> -5495
$ws.Cells.Item(137, 8).Value2 = 1834.6  # H137: 1838.68 -> 1834.6
$ws.Cells.Item(137, 9).Value2 = 1459.3846  # I137: 1467.2307 -> 1459.3846
$ws.Cells.Item(137, 11).Value2 = 4378.1538  # K137: 4401.6921 -> 4378.1538
$ws.Cells.Item(137, 13).Value2 = -1828.1538  # M137: -1851.6921 -> -1828.1538
$ws.Cells.Item(138, 8).Value2 = 2311.0454  # H138: 2303.418 -> 2311.0454
$ws.Cells.Item(138, 10).Value2 = 2257.4546  # J138: 2249.2856 -> 2257.4546
$ws.Cells.Item(138, 12).Value2 = 6772.3638  # L138: 6747.8568 -> 6772.3638
$ws.Cells.Item(138, 14).Value2 = -17052.3638  # N138: -17027.8568 -> -17052.3638
$ws.Cells.Item(141, 8).Value2 = 1080.5476  # H141: 1166.0476 -> 1080.5476
$ws.Cells.Item(141, 9).Value2 = 881.64105  # I141: 910.4054 -> 881.64105
$ws.Cells.Item(141, 10).Value2 = 3666.3333  # J141: 3057.8 -> 3666.3333
$ws.Cells.Item(141, 11).Value2 = 2644.92315  # K141: 2731.2162 -> 2644.92315
$ws.Cells.Item(141, 12).Value2 = 10998.9999  # L141: 9173.400000000001 -> 10998.9999
$ws.Cells.Item(141, 13).Value2 = 2535.07685  # M141: 2448.7838 -> 2535.07685
$ws.Cells.Item(141, 14).Value2 = -21358.9999  # N141: -19533.4 -> -21358.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value2 = 2607.25  # H45: 2643.516 -> 2607.25
$ws.Cells.Item(45, 10).Value2 = 3018.4443  # J45: 3210.375 -> 3018.4443
$ws.Cells.Item(45, 12).Value2 = 3018.4443  # L45: 3210.375 -> 3018.4443
$ws.Cells.Item(45, 14).Value2 = -3772.4443  # N45: -3964.375 -> -3772.4443
$ws.Cells.Item(76, 8).Value2 = 11874.5  # H76: 19998 -> 11874.5
$ws.Cells.Item(76, 10).Value2 = 11874.5  # J76: 19998 -> 11874.5
$ws.Cells.Item(76, 12).Value2 = 11874.5  # L76: 19998 -> 11874.5
$ws.Cells.Item(76, 14).Value2 = -12550.5  # N76: -20674 -> -12550.5
$ws.Cells.Item(79, 8).Value2 = 11874.5  # H79: 19998 -> 11874.5
$ws.Cells.Item(79, 10).Value2 = 11874.5  # J79: 19998 -> 11874.5
$ws.Cells.Item(79, 12).Value2 = 11874.5  # L79: 19998 -> 11874.5
$ws.Cells.Item(79, 14).Value2 = -14214.5  # N79: -22338 -> -14214.5
$ws.Cells.Item(122, 8).Value2 = 2054.8845  # H122: 2034.1666 -> 2054.8845
$ws.Cells.Item(122, 9).Value2 = 1889.2222  # I122: 1785.25 -> 1889.2222
$ws.Cells.Item(122, 10).Value2 = 2427.625  # J122: 2532 -> 2427.625
$ws.Cells.Item(122, 11).Value2 = 5667.6666  # K122: 5355.75 -> 5667.6666
$ws.Cells.Item(122, 12).Value2 = 7282.875  # L122: 7596 -> 7282.875
$ws.Cells.Item(122, 13).Value2 = -3217.6666  # M122: -2905.75 -> -3217.6666
$ws.Cells.Item(122, 14).Value2 = -12182.875  # N122: -12496 -> -12182.875
$ws.Cells.Item(132, 8).Value2 = 27507.861  # H132: 27508.242 -> 27507.861
$ws.Cells.Item(132, 9).Value2 = 1352.2766  # I132: 1379.5 -> 1352.2766
$ws.Cells.Item(132, 10).Value2 = 139263.55  # J132: 127668.414 -> 139263.55
$ws.Cells.Item(132, 11).Value2 = 4056.8298  # K132: 4138.5 -> 4056.8298
$ws.Cells.Item(132, 12).Value2 = 417790.65  # L132: 383005.242 -> 417790.65
$ws.Cells.Item(132, 13).Value2 = -1526.8298  # M132: -1608.5 -> -1526.8298
$ws.Cells.Item(132, 14).Value2 = -422850.65  # N132: -388065.242 -> -422850.65

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(88, 8).Value2 = 39500  # H88: 36333 -> 39500
$ws.Cells.Item(88, 10).Value2 = 39500  # J88: 36333 -> 39500
$ws.Cells.Item(88, 12).Value2 = 39500  # L88: 36333 -> 39500
$ws.Cells.Item(88, 14).Value2 = -40312  # N88: -37145 -> -40312
$ws.Cells.Item(91, 8).Value2 = 39500  # H91: 36333 -> 39500
$ws.Cells.Item(91, 10).Value2 = 39500  # J91: 36333 -> 39500
$ws.Cells.Item(91, 12).Value2 = 39500  # L91: 36333 -> 39500
$ws.Cells.Item(91, 14).Value2 = -42308  # N91: -39141 -> -42308
$ws.Cells.Item(92, 8).Value2 = 30000  # H92: 25200.5 -> 30000
$ws.Cells.Item(92, 10).Value2 = 30000  # J92: 25200.5 -> 30000
$ws.Cells.Item(92, 12).Value2 = 30000  # L92: 25200.5 -> 30000
$ws.Cells.Item(92, 14).Value2 = -34992  # N92: -30192.5 -> -34992
$ws.Cells.Item(134, 8).Value2 = 6944.625  # H134: 7140.0625 -> 6944.625
$ws.Cells.Item(134, 9).Value2 = 7365.2856  # I134: 8434.5 -> 7365.2856
$ws.Cells.Item(134, 10).Value2 = 4000  # J134: 3256.75 -> 4000
$ws.Cells.Item(134, 11).Value2 = 22095.8568  # K134: 25303.5 -> 22095.8568
$ws.Cells.Item(134, 12).Value2 = 12000  # L134: 9770.25 -> 12000
$ws.Cells.Item(134, 13).Value2 = -19560.8568  # M134: -22768.5 -> -19560.8568
$ws.Cells.Item(134, 14).Value2 = -17070  # N134: -14840.25 -> -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 15460.968  # H31: 12782.605 -> 15460.968
$ws.Cells.Item(31, 9).Value2 = 27908.334  # I31: 19321.545 -> 27908.334
$ws.Cells.Item(31, 11).Value2 = 27908.334  # K31: 19321.545 -> 27908.334
$ws.Cells.Item(31, 13).Value2 = -27613.334  # M31: -19026.545 -> -27613.334
$ws.Cells.Item(34, 8).Value2 = 15460.968  # H34: 12782.605 -> 15460.968
$ws.Cells.Item(34, 9).Value2 = 27908.334  # I34: 19321.545 -> 27908.334
$ws.Cells.Item(34, 11).Value2 = 27908.334  # K34: 19321.545 -> 27908.334
$ws.Cells.Item(34, 13).Value2 = -27706.334  # M34: -19119.545 -> -27706.334
$ws.Cells.Item(58, 8).Value2 = 32840.938  # H58: 17090.42 -> 32840.938
$ws.Cells.Item(58, 9).Value2 = 1776.4546  # I58: 1098.1305 -> 1776.4546
$ws.Cells.Item(58, 10).Value2 = 101182.8  # J58: 63068.25 -> 101182.8
$ws.Cells.Item(58, 11).Value2 = 1776.4546  # K58: 1098.1305 -> 1776.4546
$ws.Cells.Item(58, 12).Value2 = 101182.8  # L58: 63068.25 -> 101182.8
$ws.Cells.Item(58, 13).Value2 = -1573.4546  # M58: -895.1305 -> -1573.4546
$ws.Cells.Item(58, 14).Value2 = -101588.8  # N58: -63474.25 -> -101588.8
$ws.Cells.Item(74, 8).Value2 = 27629.076  # H74: 28855.334 -> 27629.076
$ws.Cells.Item(74, 10).Value2 = 27629.076  # J74: 28855.334 -> 27629.076
$ws.Cells.Item(74, 12).Value2 = 27629.076  # L74: 28855.334 -> 27629.076
$ws.Cells.Item(74, 14).Value2 = -29377.076  # N74: -30603.334 -> -29377.076
$ws.Cells.Item(77, 8).Value2 = 27629.076  # H77: 28855.334 -> 27629.076
$ws.Cells.Item(77, 10).Value2 = 27629.076  # J77: 28855.334 -> 27629.076
$ws.Cells.Item(77, 12).Value2 = 82887.228  # L77: 86566.00199999999 -> 82887.228
$ws.Cells.Item(77, 14).Value2 = -91623.228  # N77: -95302.00199999999 -> -91623.228
$ws.Cells.Item(88, 8).Value2 = 18562  # H88: 0 -> 18562
$ws.Cells.Item(88, 10).Value2 = 18562  # J88: 0 -> 18562
$ws.Cells.Item(88, 12).Value2 = 18562  # L88: 0 -> 18562
$ws.Cells.Item(88, 14).Value2 = -19374  # N88: None -> -19374
$ws.Cells.Item(91, 8).Value2 = 18562  # H91: 0 -> 18562
$ws.Cells.Item(91, 10).Value2 = 18562  # J91: 0 -> 18562
$ws.Cells.Item(91, 12).Value2 = 18562  # L91: 0 -> 18562
$ws.Cells.Item(91, 14).Value2 = -21370  # N91: None -> -21370
$ws.Cells.Item(92, 8).Value2 = 30101  # H92: 0 -> 30101
$ws.Cells.Item(92, 10).Value2 = 30101  # J92: 0 -> 30101
$ws.Cells.Item(92, 12).Value2 = 30101  # L92: 0 -> 30101
$ws.Cells.Item(92, 14).Value2 = -35093  # N92: None -> -35093
$ws.Cells.Item(99, 8).Value2 = 16132652  # H99: 15628744 -> 16132652
$ws.Cells.Item(99, 9).Value2 = 3323.182  # I99: 3446.8696 -> 3323.182
$ws.Cells.Item(99, 10).Value2 = 55559900  # J99: 55560056 -> 55559900
$ws.Cells.Item(99, 11).Value2 = 3323.182  # K99: 3446.8696 -> 3323.182
$ws.Cells.Item(99, 12).Value2 = 55559900  # L99: 55560056 -> 55559900
$ws.Cells.Item(99, 13).Value2 = -1825.182  # M99: -1948.8696 -> -1825.182
$ws.Cells.Item(99, 14).Value2 = -55562896  # N99: -55563052 -> -55562896
$ws.Cells.Item(126, 8).Value2 = 16132652  # H126: 15628744 -> 16132652
$ws.Cells.Item(126, 9).Value2 = 3323.182  # I126: 3446.8696 -> 3323.182
$ws.Cells.Item(126, 10).Value2 = 55559900  # J126: 55560056 -> 55559900
$ws.Cells.Item(126, 11).Value2 = 9969.545999999998  # K126: 10340.6088 -> 9969.545999999998
$ws.Cells.Item(126, 12).Value2 = 166679700  # L126: 166680168 -> 166679700
$ws.Cells.Item(126, 13).Value2 = -7499.545999999998  # M126: -7870.6088 -> -7499.545999999998
$ws.Cells.Item(126, 14).Value2 = -166684640  # N126: -166685108 -> -166684640
$ws.Cells.Item(136, 8).Value2 = 32840.938  # H136: 17090.42 -> 32840.938
$ws.Cells.Item(136, 9).Value2 = 1776.4546  # I136: 1098.1305 -> 1776.4546
$ws.Cells.Item(136, 10).Value2 = 101182.8  # J136: 63068.25 -> 101182.8
$ws.Cells.Item(136, 11).Value2 = 5329.3638  # K136: 3294.3915 -> 5329.3638
$ws.Cells.Item(136, 12).Value2 = 303548.4  # L136: 189204.75 -> 303548.4
$ws.Cells.Item(136, 13).Value2 = -2779.3638  # M136: -744.3914999999997 -> -2779.3638
$ws.Cells.Item(136, 14).Value2 = -308648.4  # N136: -194304.75 -> -308648.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value2 = 821.71  # H131: 135965.12 -> 821.71
$ws.Cells.Item(131, 10).Value2 = 821.71  # J131: 135965.12 -> 821.71
$ws.Cells.Item(131, 12).Value2 = 2465.13  # L131: 407895.36 -> 2465.13
$ws.Cells.Item(131, 14).Value2 = -12545.13  # N131: -417975.36 -> -12545.13

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value2 = 8925  # H70: 9472.5 -> 8925
$ws.Cells.Item(70, 10).Value2 = 4558.2  # J70: 4698.375 -> 4558.2
$ws.Cells.Item(70, 12).Value2 = 4558.2  # L70: 4698.375 -> 4558.2
$ws.Cells.Item(70, 14).Value2 = -5098.2  # N70: -5238.375 -> -5098.2
$ws.Cells.Item(73, 8).Value2 = 8925  # H73: 9472.5 -> 8925
$ws.Cells.Item(73, 10).Value2 = 4558.2  # J73: 4698.375 -> 4558.2
$ws.Cells.Item(73, 12).Value2 = 4558.2  # L73: 4698.375 -> 4558.2
$ws.Cells.Item(73, 14).Value2 = -6430.2  # N73: -6570.375 -> -6430.2
$ws.Cells.Item(80, 8).Value2 = 4322.222  # H80: 4310.6 -> 4322.222
$ws.Cells.Item(80, 10).Value2 = 4583.3335  # J80: 4529.4287 -> 4583.3335
$ws.Cells.Item(80, 12).Value2 = 4583.3335  # L80: 4529.4287 -> 4583.3335
$ws.Cells.Item(80, 14).Value2 = -6579.3335  # N80: -6525.4287 -> -6579.3335
$ws.Cells.Item(83, 8).Value2 = 4322.222  # H83: 4310.6 -> 4322.222
$ws.Cells.Item(83, 10).Value2 = 4583.3335  # J83: 4529.4287 -> 4583.3335
$ws.Cells.Item(83, 12).Value2 = 22916.6675  # L83: 22647.1435 -> 22916.6675
$ws.Cells.Item(83, 14).Value2 = -32900.6675  # N83: -32631.1435 -> -32900.6675
$ws.Cells.Item(132, 8).Value2 = 16313.973  # H132: 15906.237 -> 16313.973
$ws.Cells.Item(132, 9).Value2 = 2988.8462  # I132: 3001.1538 -> 2988.8462
$ws.Cells.Item(132, 10).Value2 = 47809.727  # J132: 43867.25 -> 47809.727
$ws.Cells.Item(132, 11).Value2 = 8966.5386  # K132: 9003.4614 -> 8966.5386
$ws.Cells.Item(132, 12).Value2 = 143429.181  # L132: 131601.75 -> 143429.181
$ws.Cells.Item(132, 13).Value2 = -6436.5386  # M132: -6473.4614 -> -6436.5386
$ws.Cells.Item(132, 14).Value2 = -148489.181  # N132: -136661.75 -> -148489.181
$ws.Cells.Item(136, 8).Value2 = 19212.5  # H136: 19389.8 -> 19212.5
$ws.Cells.Item(136, 10).Value2 = 19212.5  # J136: 19389.8 -> 19212.5
$ws.Cells.Item(136, 12).Value2 = 57637.5  # L136: 58169.39999999999 -> 57637.5
$ws.Cells.Item(136, 14).Value2 = -62737.5  # N136: -63269.39999999999 -> -62737.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(64, 8).Value2 = 23444.111  # H64: 25733.428 -> 23444.111
$ws.Cells.Item(64, 9).Value2 = 0  # I64: 10136 -> 0
$ws.Cells.Item(64, 10).Value2 = 23444.111  # J64: 28333 -> 23444.111
$ws.Cells.Item(64, 11).Value2 = 0  # K64: 10136 -> 0
$ws.Cells.Item(64, 12).Value2 = 23444.111  # L64: 28333 -> 23444.111
$ws.Cells.Item(64, 13).ClearContents()  # M64: -9911 -> (removed)
$ws.Cells.Item(64, 14).Value2 = -23894.111  # N64: -28783 -> -23894.111
$ws.Cells.Item(67, 8).Value2 = 23444.111  # H67: 25733.428 -> 23444.111
$ws.Cells.Item(67, 9).Value2 = 0  # I67: 10136 -> 0
$ws.Cells.Item(67, 10).Value2 = 23444.111  # J67: 28333 -> 23444.111
$ws.Cells.Item(67, 11).Value2 = 0  # K67: 10136 -> 0
$ws.Cells.Item(67, 12).Value2 = 23444.111  # L67: 28333 -> 23444.111
$ws.Cells.Item(67, 13).ClearContents()  # M67: -9356 -> (removed)
$ws.Cells.Item(67, 14).Value2 = -25004.111  # N67: -29893 -> -25004.111
$ws.Cells.Item(82, 8).Value2 = 2076.6924  # H82: 2214 -> 2076.6924
$ws.Cells.Item(82, 9).Value2 = 1942.5714  # I82: 2080.875 -> 1942.5714
$ws.Cells.Item(82, 11).Value2 = 1942.5714  # K82: 2080.875 -> 1942.5714
$ws.Cells.Item(82, 13).Value2 = -1581.5714  # M82: -1719.875 -> -1581.5714
$ws.Cells.Item(85, 8).Value2 = 2076.6924  # H85: 2214 -> 2076.6924
$ws.Cells.Item(85, 9).Value2 = 1942.5714  # I85: 2080.875 -> 1942.5714
$ws.Cells.Item(85, 11).Value2 = 1942.5714  # K85: 2080.875 -> 1942.5714
$ws.Cells.Item(85, 13).Value2 = -694.5714  # M85: -832.875 -> -694.5714
$ws.Cells.Item(93, 8).Value2 = 1299.6  # H93: 1723.1111 -> 1299.6
$ws.Cells.Item(93, 9).Value2 = 1277.3334  # I93: 1418 -> 1277.3334
$ws.Cells.Item(93, 10).Value2 = 1500  # J93: 2333.3333 -> 1500
$ws.Cells.Item(93, 11).Value2 = 1277.3334  # K93: 1418 -> 1277.3334
$ws.Cells.Item(93, 12).Value2 = 1500  # L93: 2333.3333 -> 1500
$ws.Cells.Item(93, 13).Value2 = -29.33339999999998  # M93: -170 -> -29.33339999999998
$ws.Cells.Item(93, 14).Value2 = -3996  # N93: -4829.3333 -> -3996
$ws.Cells.Item(132, 8).Value2 = 1504.5294  # H132: 1508.0883 -> 1504.5294
$ws.Cells.Item(132, 9).Value2 = 1219  # I132: 1223.1724 -> 1219
$ws.Cells.Item(132, 11).Value2 = 3657  # K132: 3669.5172 -> 3657
$ws.Cells.Item(132, 13).Value2 = -1127  # M132: -1139.5172 -> -1127
$ws.Cells.Item(136, 8).Value2 = 25438.477  # H136: 35257 -> 25438.477
$ws.Cells.Item(136, 9).Value2 = 34559.867  # I136: 51374.5 -> 34559.867
$ws.Cells.Item(136, 10).Value2 = 2635  # J136: 3022 -> 2635
$ws.Cells.Item(136, 11).Value2 = 103679.601  # K136: 154123.5 -> 103679.601
$ws.Cells.Item(136, 12).Value2 = 7905  # L136: 9066 -> 7905
$ws.Cells.Item(136, 13).Value2 = -101129.601  # M136: -151573.5 -> -101129.601
$ws.Cells.Item(136, 14).Value2 = -13005  # N136: -14166 -> -13005

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value2 = 3497117.5  # H107: 3497147.8 -> 3497117.5
$ws.Cells.Item(107, 9).Value2 = 709.44446  # I107: 816.75 -> 709.44446
$ws.Cells.Item(107, 10).Value2 = 11364036  # J107: 9091277 -> 11364036
$ws.Cells.Item(107, 11).Value2 = 2128.33338  # K107: 2450.25 -> 2128.33338
$ws.Cells.Item(107, 12).Value2 = 34092108  # L107: 27273831 -> 34092108
$ws.Cells.Item(107, 13).Value2 = -208.33338  # M107: -530.25 -> -208.33338
$ws.Cells.Item(107, 14).Value2 = -34095948  # N107: -27277671 -> -34095948
